# Apply the "pago excel" registration edit to the "registro" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registro")

# Row 2: update monto, moneda, add banco origen/destino, bump numero venta
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "$"
$ws.Range("G2").Value = "Mercantil"
$ws.Range("H2").Value = "BFC comecio"
$ws.Range("I2").Value = 7

# Rows 3-5: clear all the sample/demo data entirely (values + formulas),
# but keep A3:A5 / B3:B5 formatted (date / currency) and empty.
$ws.Range("C3:I5").ClearContents()
$ws.Range("A3:B5").ClearContents()

# Move the active selection from I8 to I7
$ws.Range("I7").Select()
